# Auto-generated Excel COM-interop script
# Applies the Sagittarius_Profits market-data refresh to each item sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1821.4736
$ws.Range("J17").Value = 1821.4736
$ws.Range("L17").Value = 5464.4208
$ws.Range("N17").Value = -5800.4208
$ws.Range("H70").Value = 1714
$ws.Range("I70").Value = 666
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 1998
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -1728
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 1714
$ws.Range("I73").Value = 666
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 1998
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -1062
$ws.Range("N73").Value = -9372
$ws.Range("H106").Value = 36518.4
$ws.Range("I106").Value = 36518.4
$ws.Range("K106").Value = 36518.4
$ws.Range("M106").Value = -35887.4
$ws.Range("H113").Value = 2594.2
$ws.Range("I113").Value = 1888.6
$ws.Range("K113").Value = 1888.6
$ws.Range("M113").Value = 1365.4
$ws.Range("H118").Value = 896
$ws.Range("I118").Value = 896
$ws.Range("K118").Value = 2688
$ws.Range("M118").Value = -1031
$ws.Range("H132").Value = 906.1429000000001
$ws.Range("I132").Value = 840.6667
$ws.Range("K132").Value = 2522.0001
$ws.Range("M132").Value = 7.999899999999798
$ws.Range("H137").Value = 1624.25
$ws.Range("I137").Value = 1499
$ws.Range("J137").Value = 1666
$ws.Range("K137").Value = 4497
$ws.Range("L137").Value = 4998
$ws.Range("M137").Value = -1947
$ws.Range("N137").Value = -10098
$ws.Range("H138").Value = 3846.225
$ws.Range("I138").Value = 2050.3333
$ws.Range("K138").Value = 6150.999899999999
$ws.Range("M138").Value = -1010.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2399.5715
$ws.Range("I61").Value = 2533.0833
$ws.Range("K61").Value = 2533.0833
$ws.Range("M61").Value = -2321.0833
$ws.Range("H122").Value = 2767
$ws.Range("I122").Value = 2792
$ws.Range("J122").Value = 2654.5
$ws.Range("K122").Value = 8376
$ws.Range("L122").Value = 7963.5
$ws.Range("M122").Value = -5926
$ws.Range("N122").Value = -12863.5
$ws.Range("H132").Value = 1980.0416
$ws.Range("I132").Value = 1913.9565
$ws.Range("K132").Value = 5741.8695
$ws.Range("M132").Value = -3211.8695
$ws.Range("H136").Value = 2399.5715
$ws.Range("I136").Value = 2533.0833
$ws.Range("K136").Value = 7599.249899999999
$ws.Range("M136").Value = -5049.249899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2037.5
$ws.Range("I86").Value = 2037.5
$ws.Range("K86").Value = 2037.5
$ws.Range("M86").Value = -914.5
$ws.Range("H89").Value = 2037.5
$ws.Range("I89").Value = 2037.5
$ws.Range("K89").Value = 10187.5
$ws.Range("M89").Value = -4571.5
$ws.Range("H94").Value = 735.3
$ws.Range("I94").Value = 559.1667
$ws.Range("J94").Value = 999.5
$ws.Range("K94").Value = 559.1667
$ws.Range("L94").Value = 999.5
$ws.Range("M94").Value = -108.1667
$ws.Range("N94").Value = -1901.5
$ws.Range("H105").Value = 2069.875
$ws.Range("I105").Value = 2069.875
$ws.Range("K105").Value = 2069.875
$ws.Range("M105").Value = -322.875
$ws.Range("H134").Value = 3135.2222
$ws.Range("I134").Value = 3135.2222
$ws.Range("K134").Value = 9405.6666
$ws.Range("M134").Value = -6870.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 660.3333
$ws.Range("I22").Value = 643.1667
$ws.Range("J22").Value = 694.6667
$ws.Range("K22").Value = 643.1667
$ws.Range("L22").Value = 694.6667
$ws.Range("M22").Value = -293.1667
$ws.Range("N22").Value = -1394.6667
$ws.Range("H58").Value = 2196.7693
$ws.Range("I58").Value = 1837.4286
$ws.Range("K58").Value = 1837.4286
$ws.Range("M58").Value = -1634.4286
$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -42246
$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -131232
$ws.Range("H94").Value = 143877.62
$ws.Range("J94").Value = 5340.3335
$ws.Range("L94").Value = 5340.3335
$ws.Range("N94").Value = -6242.3335
$ws.Range("H136").Value = 2196.7693
$ws.Range("I136").Value = 1837.4286
$ws.Range("K136").Value = 5512.2858
$ws.Range("M136").Value = -2962.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 15500
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 51458.453
$ws.Range("J15").Value = 51458.453
$ws.Range("L15").Value = 51458.453
$ws.Range("N15").Value = -52034.453
$ws.Range("H81").Value = 51458.453
$ws.Range("J81").Value = 51458.453
$ws.Range("L81").Value = 51458.453
$ws.Range("N81").Value = -53454.453
$ws.Range("H84").Value = 51458.453
$ws.Range("J84").Value = 51458.453
$ws.Range("L84").Value = 154375.359
$ws.Range("N84").Value = -164359.359
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2654
$ws.Range("I46").Value = 1721.5
$ws.Range("K46").Value = 1721.5
$ws.Range("M46").Value = -1533.5
$ws.Range("H55").Value = 498.9375
$ws.Range("I55").Value = 504
$ws.Range("J55").Value = 487.8
$ws.Range("K55").Value = 504
$ws.Range("L55").Value = 487.8
$ws.Range("M55").Value = -331
$ws.Range("N55").Value = -833.8
$ws.Range("H122").Value = 5768.8076
$ws.Range("I122").Value = 4874.75
$ws.Range("J122").Value = 7199.3
$ws.Range("K122").Value = 14624.25
$ws.Range("L122").Value = 21597.9
$ws.Range("M122").Value = -12174.25
$ws.Range("N122").Value = -26497.9
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 7298.1
$ws.Range("I132").Value = 8782.143
$ws.Range("K132").Value = 26346.429
$ws.Range("M132").Value = -23816.429
$ws.Range("H136").Value = 1843
$ws.Range("I136").Value = 1106.3572
$ws.Range("K136").Value = 3319.0716
$ws.Range("M136").Value = -769.0715999999998
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 90095
$ws.Range("I70").Value = 90095
$ws.Range("K70").Value = 90095
$ws.Range("M70").Value = -89780
$ws.Range("H73").Value = 90095
$ws.Range("I73").Value = 90095
$ws.Range("K73").Value = 90095
$ws.Range("M73").Value = -89003
$ws.Range("H75").Value = 74376
$ws.Range("I75").Value = 74509
$ws.Range("K75").Value = 74509
$ws.Range("M75").Value = -73573
$ws.Range("H78").Value = 74376
$ws.Range("I78").Value = 74509
$ws.Range("K78").Value = 223527
$ws.Range("M78").Value = -218847
$ws.Range("H100").Value = 20001482
$ws.Range("I100").Value = 25001298
$ws.Range("J100").Value = 2222
$ws.Range("K100").Value = 50002596
$ws.Range("L100").Value = 4444
$ws.Range("M100").Value = -50002055
$ws.Range("N100").Value = -5526
$ws.Range("H122").Value = 2188.8
$ws.Range("I122").Value = 2188.8
$ws.Range("K122").Value = 6566.400000000001
$ws.Range("M122").Value = -4116.400000000001
$ws.Range("H136").Value = 2071.65
$ws.Range("I136").Value = 2063.9412
$ws.Range("K136").Value = 6191.823600000001
$ws.Range("M136").Value = -3641.823600000001

